$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original plain-text storage
# (values like "1.002" or "6.600" would otherwise be auto-coerced to
# numbers by Excel and lose trailing zeros / exact formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.300.24'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').Value = '1.628.59'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '297.97'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('D8').Value = '50.04'
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').Value = '0.3472'
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('D10').Value = '0.08025'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('D11').Value = '1.198'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').Value = '6.284'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').Value = '7.210'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = '0.00001188'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '1.628.22'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('D19').Value = '0.06943'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').Value = '6.600'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '12.34'
$ws.Range('E23').Value = '  -3.78%  '
$ws.Range('D24').Value = '23.302.85'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = '2.409'
$ws.Range('E25').Value = '  -3.81%  '
$ws.Range('D26').Value = '2.952'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').Value = '20.83'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('D28').Value = '150.65'
$ws.Range('E28').Value = '  -1.98%  '
$ws.Range('D29').Value = '5.141'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('D30').Value = '131.27'
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('D31').Value = '1.805.45'
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('D32').Value = '6.724'
$ws.Range('E32').Value = '  -5.13%  '
$ws.Range('D33').Value = '2.117'
$ws.Range('E33').Value = '  -5.68%  '
$ws.Range('D34').Value = '11.15'
$ws.Range('E34').Value = '  -7.86%  '
$ws.Range('D35').Value = '0.9753'
$ws.Range('E35').Value = '  -7.46%  '
$ws.Range('D36').Value = '0.02653'
$ws.Range('E36').Value = '  -5.72%  '
$ws.Range('D37').Value = '0.08737'
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('D38').Value = '0.2410'
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('D39').Value = '5.820'
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('D40').Value = '0.06691'
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('D41').Value = '12.64'
$ws.Range('E41').Value = '  -2.77%  '
$ws.Range('D42').Value = '0.6778'
$ws.Range('E42').Value = '  -3.22%  '
$ws.Range('D43').Value = '1.290'
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('D44').Value = '15.29'
$ws.Range('E44').Value = '  -4.72%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '0.6294'
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('D47').Value = '2.225'
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').Value = '3.882'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D49').Value = '0.07612'
$ws.Range('E49').Value = '  -3.84%  '
$ws.Range('D50').Value = '125.91'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').Value = '  +1.72%  '
